# Mexico Liga de Expansion - base update (03-03-2024 00:35)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 91 and 92 got swapped (the match that used to be listed on row 91 is
# now listed on row 92 and vice-versa). Column A (the running index) stays
# untouched; every other column (B..AC) is exchanged between the two rows.
# ---------------------------------------------------------------------------

# New row 91 (was row 92)
$ws.Cells.Item(91, 2).Value2 = 6924568
$ws.Cells.Item(91, 6).Value2 = "Atletico Morelia"
$ws.Cells.Item(91, 7).Value2 = "Atlante"
$ws.Cells.Item(91, 8).Value2 = 0
$ws.Cells.Item(91, 10).Value2 = "A"
$ws.Cells.Item(91, 11).Value2 = 2.4
$ws.Cells.Item(91, 12).Value2 = 3
$ws.Cells.Item(91, 13).Value2 = 2.875
$ws.Cells.Item(91, 14).Value2 = 2.7
$ws.Cells.Item(91, 15).Value2 = 3.1
$ws.Cells.Item(91, 16).Value2 = 2.8
$ws.Cells.Item(91, 17).Value2 = 0
$ws.Cells.Item(91, 18).Value2 = 1.85
$ws.Cells.Item(91, 19).Value2 = 1.95
$ws.Cells.Item(91, 20).Value2 = 2.25
$ws.Cells.Item(91, 21).Value2 = 1.975
$ws.Cells.Item(91, 22).Value2 = 1.725
$ws.Cells.Item(91, 23).Value2 = -1
$ws.Cells.Item(91, 25).Value2 = 1.8
$ws.Cells.Item(91, 26).Value2 = -1
$ws.Cells.Item(91, 27).Value2 = 0.95
$ws.Cells.Item(91, 28).Value2 = -1
$ws.Cells.Item(91, 29).Value2 = 0.7250000000000001

# New row 92 (was row 91)
$ws.Cells.Item(92, 2).Value2 = 6924569
$ws.Cells.Item(92, 6).Value2 = "Venados FC"
$ws.Cells.Item(92, 7).Value2 = "Dorados"
$ws.Cells.Item(92, 8).Value2 = 4
$ws.Cells.Item(92, 10).Value2 = "H"
$ws.Cells.Item(92, 11).Value2 = 1.615
$ws.Cells.Item(92, 12).Value2 = 4
$ws.Cells.Item(92, 13).Value2 = 4.5
$ws.Cells.Item(92, 14).Value2 = 1.5
$ws.Cells.Item(92, 15).Value2 = 4.75
$ws.Cells.Item(92, 16).Value2 = 5.75
$ws.Cells.Item(92, 17).Value2 = -1.25
$ws.Cells.Item(92, 18).Value2 = 1.925
$ws.Cells.Item(92, 19).Value2 = 1.875
$ws.Cells.Item(92, 20).Value2 = 3
$ws.Cells.Item(92, 21).Value2 = 1.75
$ws.Cells.Item(92, 22).Value2 = 1.95
$ws.Cells.Item(92, 23).Value2 = 0.5
$ws.Cells.Item(92, 25).Value2 = -1
$ws.Cells.Item(92, 26).Value2 = 0.925
$ws.Cells.Item(92, 27).Value2 = -1
$ws.Cells.Item(92, 28).Value2 = 0.75
$ws.Cells.Item(92, 29).Value2 = -1

# ---------------------------------------------------------------------------
# Row 179: odds / id / date / teams refreshed with the closing values for the
# fixture (the match has not been played yet, so FTHG/FTAG/FTR stay blank).
# ---------------------------------------------------------------------------
$ws.Cells.Item(179, 2).Value2 = 7640647
$ws.Cells.Item(179, 5).Value2 = 45356.92013888889
$ws.Cells.Item(179, 6).Value2 = "Atlante"
$ws.Cells.Item(179, 7).Value2 = "Tlaxcala FC"
$ws.Cells.Item(179, 11).Value2 = 1.4
$ws.Cells.Item(179, 12).Value2 = 4.5
$ws.Cells.Item(179, 13).Value2 = 5.75
$ws.Cells.Item(179, 14).Value2 = 1.4
$ws.Cells.Item(179, 15).Value2 = 4.333
$ws.Cells.Item(179, 16).Value2 = 6
$ws.Cells.Item(179, 17).Value2 = -1.25
$ws.Cells.Item(179, 18).Value2 = 1.875
$ws.Cells.Item(179, 19).Value2 = 1.925
$ws.Cells.Item(179, 21).Value2 = 1.85
$ws.Cells.Item(179, 22).Value2 = 1.95

# ---------------------------------------------------------------------------
# New fixtures appended as rows 180-183.
# ---------------------------------------------------------------------------

function Add-Fixture($Row, $Id, $MatchId, $Date, $Home, $Away, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $U, $V) {
    $ws.Cells.Item($Row, 1).Value2 = $Id
    $ws.Range("A91").Copy()
    $ws.Cells.Item($Row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($Row, 2).Value2 = $MatchId

    $ws.Cells.Item($Row, 3).Value2 = "Mexico Liga de Expansion"
    $ws.Cells.Item($Row, 4).Value2 = "Mexico Liga de Expansion"

    $ws.Cells.Item($Row, 5).Value2 = $Date
    $ws.Range("E91").Copy()
    $ws.Cells.Item($Row, 5).PasteSpecial(-4122)

    $ws.Cells.Item($Row, 6).Value2 = $Home
    $ws.Cells.Item($Row, 7).Value2 = $Away

    $ws.Cells.Item($Row, 11).Value2 = $K
    $ws.Cells.Item($Row, 12).Value2 = $L
    $ws.Cells.Item($Row, 13).Value2 = $M
    $ws.Cells.Item($Row, 14).Value2 = $N
    $ws.Cells.Item($Row, 15).Value2 = $O
    $ws.Cells.Item($Row, 16).Value2 = $P
    $ws.Cells.Item($Row, 17).Value2 = $Q
    $ws.Cells.Item($Row, 18).Value2 = $R
    $ws.Cells.Item($Row, 19).Value2 = $S
    $ws.Cells.Item($Row, 20).Value2 = 2.5
    $ws.Cells.Item($Row, 21).Value2 = $U
    $ws.Cells.Item($Row, 22).Value2 = $V
    $ws.Cells.Item($Row, 23).Value2 = 0
    $ws.Cells.Item($Row, 24).Value2 = 0
    $ws.Cells.Item($Row, 25).Value2 = 0
    $ws.Cells.Item($Row, 26).Value2 = 0
    $ws.Cells.Item($Row, 27).Value2 = 0
}

Add-Fixture 180 178 7641687 45357.00347222222 "Club Celaya" "Venados FC" `
    1.65 3.5 4.5 1.65 3.5 4.5 -0.75 1.825 1.975 1.975 1.825

Add-Fixture 181 179 7641688 45357.92013888889 "Correcaminos" "Tepatitlan FC" `
    1.571 3.6 5.25 1.571 3.6 5.25 -1 2 1.8 1.975 1.825

Add-Fixture 182 180 7641689 45358.00347222222 "Dorados" "Atletico Morelia" `
    3.6 3.5 1.833 3.6 3.5 1.833 0.5 1.925 1.875 1.95 1.85

Add-Fixture 183 181 7641691 45359.00347222222 "Tapatio" "Oaxaca" `
    1.571 3.75 4.75 1.571 3.75 4.75 -1 2 1.8 1.8 2
